$d = $word.ActiveDocument

# --- Step 1: remove the _GoBack bookmark from the "push origin master" paragraph ---
$pushPara = $d.Paragraphs.Item(38)
$pushXml = @'
<w:p w:rsidR="00300E1E" w:rsidRDefault="00300E1E" w:rsidP="00300E1E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> push origin master</w:t></w:r></w:p>
'@
$pushPara.Range.InsertXML($pushXml)

# --- Step 2: replace the trailing empty paragraph with the new block of paragraphs ---
$lastPara = $d.Paragraphs.Last
$newXml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">To create new branch </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> branch &lt;branch-name&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>To switch to created branch</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> checkout &lt;branch-name&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Create &amp; immediately switch to new branch</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> checkout -b &lt;branch-name&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:color w:val="24292F"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>pushing a new branch, this creates an upstream tracking branch with a lasting relationship to your local branch</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="24292F"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="24292F"/>
    </w:rPr>
    <w:t>git push -u origin [branch</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="24292F"/>
    </w:rPr>
    <w:t>-name</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="24292F"/>
    </w:rPr>
    <w:t>]</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="24292F"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>

'@
$lastPara.Range.InsertXML($newXml)

"done"
